$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateNumberFormat = $ws.Range("A5885").NumberFormat

$data = @"
5886|44176|Bánovce nad Bebravou|20
5887|44176|Banská Bystrica|106
5888|44176|Banská Štiavnica|13
5889|44176|Bardejov|44
5890|44176|Bratislava|213
5891|44176|Brezno|31
5892|44176|Bytča|13
5893|44176|Čadca|75
5894|44176|Detva|12
5895|44176|Dolný Kubín|7
5896|44176|Dunajská Streda|64
5897|44176|Galanta|53
5898|44176|Gelnica|11
5899|44176|Hlohovec|32
5900|44176|Humenné|48
5901|44176|Ilava|67
5902|44176|Kežmarok|26
5903|44176|Komárno|55
5904|44176|Košice|99
5905|44176|Košice - okolie|57
5906|44176|Krupina|17
5907|44176|Kysucké Nové Mesto|12
5908|44176|Levice|42
5909|44176|Levoča|10
5910|44176|Liptovský Mikuláš|65
5911|44176|Lučenec|28
5912|44176|Malacky|34
5913|44176|Martin|86
5914|44176|Medzilaborce|7
5915|44176|Michalovce|46
5916|44176|Myjava|57
5917|44176|Námestovo|7
5918|44176|Nitra|164
5919|44176|Nové Mesto nad Váhom|134
5920|44176|Nové Zámky|60
5921|44176|Partizánske|13
5922|44176|Pezinok|18
5923|44176|Piešťany|32
5924|44176|Poltár|7
5925|44176|Poprad|68
5926|44176|Považská Bystrica|69
5927|44176|Prešov|158
5928|44176|Prievidza|151
5929|44176|Púchov|84
5930|44176|Revúca|16
5931|44176|Rimavská Sobota|62
5932|44176|Rožňava|16
5933|44176|Ružomberok|34
5934|44176|Sabinov|16
5935|44176|Senec|61
5936|44176|Senica|67
5937|44176|Skalica|22
5938|44176|Snina|23
5939|44176|Sobrance|7
5940|44176|Spišská Nová Ves|67
5941|44176|Stará Ľubovňa|27
5942|44176|Stropkov|26
5943|44176|Svidník|48
5944|44176|Šaľa|22
5945|44176|Topoľčany|35
5946|44176|Trebišov|52
5947|44176|Trenčín|155
5948|44176|Trnava|114
5949|44176|Turčianske Teplice|7
5950|44176|Tvrdošín|6
5951|44176|Veľký Krtíš|33
5952|44176|Vranov nad Topľou|59
5953|44176|Zlaté Moravce|31
5954|44176|Zvolen|73
5955|44176|Žarnovica|8
5956|44176|Žiar nad Hronom|62
5957|44176|Žilina|143
"@

$lines = $data -split "`n"

foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    $rowNum = [int]$parts[0]
    $dateVal = [int]$parts[1]
    $district = $parts[2]
    $count = [int]$parts[3]

    $ws.Cells.Item($rowNum, 1).Value = $dateVal
    $ws.Cells.Item($rowNum, 1).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($rowNum, 2).Value = $district
    $ws.Cells.Item($rowNum, 3).Value = $count
}
